# Tratando termos banidos, termos de busca, precos e ajuste na base de dados
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Ajuste no termo de busca: "iphone 12 64gb" -> "iphone 12 64 gb"
$ws.Range("A2").Value = "iphone 12 64 gb"

# Atualiza a celula/selecao ativa conforme a ultima edicao
$ws.Activate()
$ws.Range("D8").Select()
